$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "actual time" (column D) for entries that have been completed so far.
$ws.Range("D3").Value = 0.5
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("D9").Value = 1

# Reflect where the user left the selection while editing.
$ws.Range("D8").Select()
